$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate "Test Case for Bootstrap" to create the new
#    "Test Case for Student breakdown" sheet (placed right after Bootstrap),
#    then add a fresh blank sheet "Test Case for Agd" after that.
# ---------------------------------------------------------------------------
$bootstrap = $wb.Worksheets.Item("Test Case for Bootstrap")
$bootstrap.Copy([System.Reflection.Missing]::Value, $bootstrap)
$breakdown = $wb.Worksheets.Item($bootstrap.Index + 1)
$breakdown.Name = "Test Case for Student breakdown"

$agd = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $breakdown)
$agd.Name = "Test Case for Agd"

# ---------------------------------------------------------------------------
# 2. Populate the new "Test Case for Student breakdown" sheet with the new
#    test cases, and clear out the inherited Bootstrap content/results.
# ---------------------------------------------------------------------------

# New test-case rows (time/date validation test cases).
$breakdown.Range("B3").Value2 = "User: Correct time, Correct date"
$breakdown.Range("C3").Value2 = "time: 10:00:02 AM`ndate: 06-02-2017"
$breakdown.Range("D3").Value2 = "Enter time into time field and enter date under date field "
$breakdown.Range("E3").Value2 = "User able to view four different basic location reports"
$breakdown.Range("F3").Value2 = "User able to view four different basic location reports"

$breakdown.Range("B4").Value2 = "User: Wrong time, Correct date"
$breakdown.Range("C4").Value2 = "time: 10:00 AM`ndate: 06-02-2018"
$breakdown.Range("D4").Value2 = "Enter time into time field and enter date under date field "
$breakdown.Range("E4").Value2 = "User able to view four different basic location reports"
$breakdown.Range("F4").Value2 = "User unable to view four different basic location reports"

$breakdown.Range("B5").Value2 = "User: Correct time, Wrong date"
$breakdown.Range("C5").Value2 = "time: 10:00:02 AM`ndate: 2019-02-06"
$breakdown.Range("D5").Value2 = "Enter time into time field and enter date under date field "
$breakdown.Range("E5").Value2 = "User able to view four different basic location reports"
$breakdown.Range("F5").Value2 = "User unable to view four different basic location reports"

# These are brand-new test cases -- no Pass/Fail results recorded yet, and
# the remainder of the sheet (rows 6-32, carried over from Bootstrap) is
# cleared out so it is ready for future entries.
$breakdown.Range("G3:K32").ClearContents()
$breakdown.Range("B6:F32").ClearContents()

# Row heights for the three wrapped test-case rows.
$breakdown.Rows.Item(3).RowHeight = 72.5
$breakdown.Rows.Item(4).RowHeight = 72.5
$breakdown.Rows.Item(5).RowHeight = 72.5

# Column B needs to be a little wider to fit the new descriptions.
$breakdown.Columns.Item(2).ColumnWidth = 28.1796875

# ---------------------------------------------------------------------------
# 3. Update selections / active-tab state to reflect where the editor's
#    cursor ended up on each sheet.
# ---------------------------------------------------------------------------
$login = $wb.Worksheets.Item("Test Case for Login Validation")
$login.Activate() | Out-Null
$login.Range("C3").Select() | Out-Null

$bootstrap.Activate() | Out-Null
$bootstrap.Range("C4").Select() | Out-Null

$breakdown.Activate() | Out-Null
$breakdown.Range("B6").Select() | Out-Null
